$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$months11 = @("Jan","Feb","Mar","Apr","May","Jun","Jul","Aug","Sep","Oct","Nov")

function Set-MonthHeaderRow($row) {
    # Clone formatting from the existing header row (B1:M1) which already
    # carries style index 1 (bold, bordered, centered).
    $ws.Range("B1:M1").Copy()
    $ws.Range("B" + $row + ":M" + $row).PasteSpecial(-4122)

    for ($i = 0; $i -lt $months11.Length; $i++) {
        $col = 2 + $i   # B = 2
        $ws.Cells.Item($row, $col).Value = $months11[$i]
    }
    $ws.Cells.Item($row, 13).Value = "avg"  # M column
}

function Set-DataRow($row, $label, $values) {
    if ($label -ne $null) {
        # Clone formatting from an existing label cell (A2) with style index 1.
        $ws.Range("A2").Copy()
        $ws.Cells.Item($row, 1).PasteSpecial(-4122)
        $ws.Cells.Item($row, 1).Value = $label
    }
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$excel.CutCopyMode = 0

# --- Block 1: rows 13-17 ---
Set-MonthHeaderRow 13
Set-DataRow 14 "Verbal Excellence" @(69,92,73,75,55,93,85,82,62,75,75,77)
Set-DataRow 15 "Avoid Rude/Unprofessional Behavior/Approach (ARU)" @(100,100,100,100,100,93,92,100,100,100,100,98)
$ws.Cells.Item(16, 13).Value = 87.5
$ws.Cells.Item(17, 13).Value = 4.38

# --- Block 2: rows 20-24 ---
Set-MonthHeaderRow 20
Set-DataRow 21 "Courteous Approach" @(69,92,100,92,82,93,85,91,69,75,67,85)
Set-DataRow 22 "Active Listening" @(100,100,100,100,100,93,100,100,100,92,100,98)
$ws.Cells.Item(23, 13).Value = 91.5
$ws.Cells.Item(24, 13).Value = 4.58

# --- Block 3: rows 27-31 ---
Set-MonthHeaderRow 27
Set-DataRow 28 "Correct & Complete Information For Resolution (CCIR)" @(100,100,100,100,100,100,100,100,100,92,100,99)
Set-DataRow 29 "Identification and Action for Resolution" @(100,100,100,100,100,100,100,100,100,100,100,100)
$ws.Cells.Item(30, 13).Value = 99.5
$ws.Cells.Item(31, 13).Value = 4.97

# --- Block 4: rows 34-37 ---
Set-MonthHeaderRow 34
Set-DataRow 35 "Ownership & Proctiveness (OP)" @(100,100,93,100,100,93,100,82,100,92,83,94)
$ws.Cells.Item(36, 13).Value = 94
$ws.Cells.Item(37, 13).Value = 4.7

$excel.CutCopyMode = 0

# Remove December column (M) from the first table (rows 1-8) -- done last so
# the earlier header-row copies (which source from row 1) still see the
# original M1 "Dec" style/formatting before it is cleared.
$ws.Range("M1:M8").Clear()
